# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for both the zh-cn and de-de language sheets to reflect the new
# report generation timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-22 14:02:31"
$wsZhCn.Range("G2").Value = "2016-02-22 14:03:21"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-22 14:02:43"
$wsDeDe.Range("G2").Value = "2016-02-22 14:03:42"
